$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 2.1
    "I2" = 4
    "L2" = 4.75
    "M2" = 1.11
    "N2" = 6.5
    "AA2" = 2.25
    "AB2" = 1.57
    "AD2" = 8.5
    "AF2" = 19
    "AK2" = 21
    "AN2" = 8.5
    "G3" = 2.2
    "H3" = 3.1
    "I3" = 3.4
    "J3" = 3.1
    "L3" = 4.33
    "M3" = 1.11
    "N3" = 6.5
    "AA3" = 2.1
    "AB3" = 1.67
    "AD3" = 9.5
    "AF3" = 21
    "AJ3" = 6
    "AL3" = 67
    "AO3" = 15
    "G4" = 2.7
    "I4" = 2.9
    "J4" = 3.75
    "AI4" = 5.5
    "AO4" = 12
    "G7" = 2.4
    "H7" = 2.75
    "I7" = 3.4
    "M7" = 1.17
    "N7" = 5
    "S7" = 3.4
    "T7" = 1.33
    "W7" = 7
    "X7" = 1.1
    "AA7" = 2.63
    "AB7" = 1.44
    "AC7" = 5
    "AD7" = 9.5
    "AE7" = 12
    "AF7" = 23
    "AH7" = 51
    "AI7" = 4.75
    "AK7" = 23
    "AO7" = 15
    "AP7" = 15
    "AR7" = 41
    "G8" = 2.63
    "H8" = 3.2
    "I8" = 2.63
    "J8" = 3.5
    "K8" = 2
    "L8" = 3.5
    "M8" = 1.08
    "N8" = 8
    "O8" = 1.4
    "P8" = 2.75
    "S8" = 2.35
    "T8" = 1.57
    "W8" = 4.33
    "X8" = 1.2
    "Y8" = 1.5
    "Z8" = 2.5
    "AA8" = 1.95
    "AB8" = 1.8
    "AC8" = 7.5
    "AE8" = 11
    "AI8" = 8
    "AJ8" = 6
    "AK8" = 17
    "AM8" = 401
    "AN8" = 7.5
    "AP8" = 11
    "G9" = 2.35
    "I9" = 3
    "Q9" = 1.83
    "R9" = 2.03
    "W9" = 4.5
    "X9" = 1.18
    "AA9" = 2.1
    "AB9" = 1.67
    "AN9" = 7.5
    "G10" = 3.4
    "I10" = 2.2
    "J10" = 4.33
    "L10" = 3
    "Q10" = 1.9
    "R10" = 1.95
    "AA10" = 2.1
    "AB10" = 1.67
    "AC10" = 7.5
    "AD10" = 15
    "AE10" = 13
    "AH10" = 41
    "AI10" = 7
    "AJ10" = 6
    "AK10" = 19
    "AL10" = 67
    "AN10" = 6
    "AO10" = 9.5
    "AP10" = 10
    "AQ10" = 21
    "G11" = 1.8
    "I11" = 4.75
    "J11" = 2.6
    "L11" = 5.5
    "AA11" = 2.2
    "AB11" = 1.62
    "AG11" = 17
    "AH11" = 41
    "AK11" = 21
    "AL11" = 81
    "AO11" = 23
    "H12" = 3.3
    "I12" = 3.2
    "K12" = 2.05
    "M12" = 1.07
    "N12" = 8.5
    "O12" = 1.36
    "P12" = 3
    "S12" = 2.2
    "T12" = 1.65
    "W12" = 4
    "X12" = 1.22
    "Y12" = 1.5
    "Z12" = 2.5
    "AA12" = 1.95
    "AB12" = 1.8
    "AC12" = 7
    "AD12" = 10
    "AG12" = 19
    "AI12" = 8.5
    "AJ12" = 6.5
    "AK12" = 17
    "AL12" = 51
    "AM12" = 351
    "AN12" = 8.5
    "AQ12" = 34
    "O14" = 1.36
    "P14" = 3
    "AM14" = 800
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
